$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update movie names for existing rows
$ws.Range("B2").Value = "daredevil"
$ws.Range("B3").Value = "superman"
$ws.Range("B4").Value = "justice league"

# Remove the now-unused rows 5-7 (clear contents so the used range shrinks)
$ws.Range("A5:B7").ClearContents()

# Move selection to reflect the new extent of the list
$ws.Range("B5").Select()
